# Rename the three picture shapes (Pearson/Edexcel logo in both footers,
# BTec logo in the first-page header) as captured by the commit:
#   footer (first page)   id=3  image2.png -> image1.png
#   footer (default/odd)  id=2  image2.png -> image1.png
#   header (first page)   id=1  image1.jpg -> image2.jpg
#
# InlineShape has no settable Name in the Word object model, so each
# picture is round-tripped through Shape (which does expose .Name) and
# then converted back to an inline shape so the <wp:inline> layout is
# preserved exactly as it was.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlinePicture($range, $newName) {
    $inline = $range.InlineShapes.Item(1)
    $shape = $inline.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape() | Out-Null
}

# Default footer (footer2.xml, docPr id="2")
Rename-InlinePicture $sec.Footers.Item(1).Range "image1.png"

# First-page footer (footer1.xml, docPr id="3")
Rename-InlinePicture $sec.Footers.Item(2).Range "image1.png"

# First-page header (header1.xml, docPr id="1")
Rename-InlinePicture $sec.Headers.Item(2).Range "image2.jpg"
